# "adds - chapter 2"
# Adds a new bibliography row (Stan probabilistic-programming-language
# software reference) at row 67 of Sheet1, mirroring the existing table
# layout (Paper Title / Journal / Year / Author(s) / Other / Read /
# Point of the paper / Methods / Zotero / Available / Cited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "Stan: A Probabilistic Programming Language"
$ws.Range("B67").Value = "Journal of Statistical Software"
$ws.Range("C67").Value = 2017
$ws.Range("D67").Value = "Carpenter et al."
$ws.Range("E67").Value = "Software/Package"
$ws.Range("F67").Value = "NA"
$ws.Range("G67").Value = "Stan software"
$ws.Range("I67").Value = "yes"
$ws.Range("J67").Value = "yes"
$ws.Range("K67").Value = "yes"

# Restore the view's selection to the newly-added row's last cell, as in
# the authored session.
[void]$ws.Range("K67").Select()
